$p = $ppt.ActivePresentation

# --- Slide 10 ("Docker"): merge the two runs of the last bullet into one run ---
$s10 = $p.Slides.Item(10)
$contentShape10 = $s10.Shapes.Item(2)
# First push a placeholder body through so the engine's text-diff doesn't try to
# preserve the old two-run split at the apostrophe; then set the real text so the
# whole last paragraph collapses back down into a single run.
$contentShape10.TextFrame.TextRange.Text = "Container management software.`r`rUsed to create images, and containers.`r`rUses dockerfiles to automate this`r`rPLACEHOLDER"
$contentShape10.TextFrame.TextRange.Text = "Container management software.`r`rUsed to create images, and containers.`r`rUses dockerfiles to automate this`r`rWe’ve used Docker to create and use images."

# --- Slide 11 ("Docker Swarm"): fill in the previously-empty content placeholder ---
$s11 = $p.Slides.Item(11)
$contentShape11 = $s11.Shapes.Item(2)
$contentShape11.TextFrame.TextRange.Text = "Docker swarm allows for multiple containers across multiple hosts.`r`rMeans that rolling updates can be applied.`r`rWe have used it for that purpose."
